$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row text updates (row 1) ---
# "X to gw quality" headers become short labels; "large growers" columns become
# "investor growers" columns.
$ws.Range("B1").Value = "rural communities"
$ws.Range("C1").Value = "small growers"
$ws.Range("D1").Value = "investor growers"
$ws.Range("E1").Value = "small growers (white area)"
$ws.Range("F1").Value = "investor growers (white area)"
$ws.Range("G1").Value = "municipalities"
# H1 text ("other dischargers") is unchanged; only its format changes below.

# --- row 1 height shrinks now that the header text is shorter ---
$ws.Rows.Item(1).RowHeight = 58.5

# --- header font: switch C1:H1 to an explicit (non theme-scheme) Calibri 11 ---
foreach ($addr in @("C1", "D1", "E1", "F1", "G1", "H1")) {
    $f = $ws.Range($addr).Font
    $f.Name = "Calibri"
    $f.Size = 11
}

# --- header borders ---
# C1: right edge medium -> thin; top edge medium -> none; bottom edge black -> grey
$c1 = $ws.Range("C1")
$c1.Borders.Item(10).Weight = 2
$c1.Borders.Item(10).Color = 0
$c1.Borders.Item(8).LineStyle = -4142
$c1.Borders.Item(9).Color = 13421772

# D1,E1,F1,G1: right edge medium -> thin; bottom edge black -> grey (top/left unchanged)
foreach ($addr in @("D1", "E1", "F1", "G1")) {
    $r = $ws.Range($addr)
    $r.Borders.Item(10).Weight = 2
    $r.Borders.Item(10).Color = 0
    $r.Borders.Item(9).Color = 13421772
}

# H1: right edge medium -> thin; top/bottom go from none -> medium grey
$h1 = $ws.Range("H1")
$h1.Borders.Item(10).Weight = 2
$h1.Borders.Item(10).Color = 0
$h1.Borders.Item(8).LineStyle = 1
$h1.Borders.Item(8).Weight = -4138
$h1.Borders.Item(8).Color = 13421772
$h1.Borders.Item(9).LineStyle = 1
$h1.Borders.Item(9).Weight = -4138
$h1.Borders.Item(9).Color = 13421772

# --- restore the pane/selection state (frozen header, scrolled near the top) ---
$ws.Range("M5").Select()
